# Deploy the implementation guide.
#
# Renames the "Include" sheet, refreshes the Metadata sheet's Date/Contact
# values, and inserts a new "Jurisdiction" row into the Metadata table.

$wb = $excel.ActiveWorkbook

# --- Rename the second worksheet ("Include from Ferlab.bio CodeS" -> "Include #0") ---
$wsInclude = $wb.Worksheets.Item(2)
$wsInclude.Name = "Include #0"

# --- Update the Metadata sheet ---
$wsMeta = $wb.Worksheets.Item(1)

# Date value (row 8, column B)
$wsMeta.Range("B8").Value = "2024-10-02T15:04:17+00:00"

# Contact value (row 10, column B)
$wsMeta.Range("B10").Value = "Ferlab.bio (http://example.org/example-publisher)"

# Insert a new row before the "Description" row (currently row 11) for "Jurisdiction"
$wsMeta.Cells.Item(11, 1).EntireRow.Insert()

$wsMeta.Range("A11").Value = "Jurisdiction"
$wsMeta.Range("B11").Value = ""

# Copy the row formatting (border/alignment/style) from the row below, which
# holds what used to be row 11 before the insert, so the new row matches the
# rest of the table's look.
$wsMeta.Range("A12:B12").Copy()
$wsMeta.Range("A11:B11").PasteSpecial(-4122)
